$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "26.714.50"
$ws.Range("E2").Value = "  +0.25%  "
Set-TextValue $ws.Range("D3") "1.641.25"
$ws.Range("E3").Value = "  -0.10%  "
$ws.Range("E4").Value = "  +0.25%  "
Set-TextValue $ws.Range("D5") "216.36"
$ws.Range("E5").Value = "  +0.66%  "
$ws.Range("E6").Value = "  -0.78%  "
$ws.Range("E7").Value = "  +0.18%  "
$ws.Range("E8").Value = "  -0.36%  "
$ws.Range("E9").Value = "  -0.06%  "
$ws.Range("E10").Value = "  -0.85%  "
Set-TextValue $ws.Range("D11") "0.0842"
$ws.Range("E11").Value = "  +0.06%  "
Set-TextValue $ws.Range("D12") "1.642.47"
$ws.Range("E12").Value = "  +0.45%  "
$ws.Range("E13").Value = "  -0.87%  "
$ws.Range("E14").Value = "  -0.72%  "
Set-TextValue $ws.Range("D15") "64.52"
$ws.Range("E15").Value = "  -1.44%  "
Set-TextValue $ws.Range("D16") "26.712.27"
$ws.Range("E16").Value = "  +0.08%  "
Set-TextValue $ws.Range("D18") "213.91"
$ws.Range("E18").Value = "  -0.98%  "
$ws.Range("E19").Value = "  +0.13%  "
$ws.Range("E20").Value = "  +0.21%  "
Set-TextValue $ws.Range("D21") "2.46"
$ws.Range("E21").Value = "  +13.05%  "
Set-TextValue $ws.Range("D22") "6.24"
$ws.Range("E22").Value = "  -0.88%  "
$ws.Range("E23").Value = "  -2.34%  "
Set-TextValue $ws.Range("D24") "145.78"
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("E25").Value = "  +0.36%  "
Set-TextValue $ws.Range("D26") "0.118"
$ws.Range("E26").Value = "  -1.73%  "
$ws.Range("E27").Value = "  -0.68%  "
$ws.Range("E28").Value = "  -0.92%  "
$ws.Range("E29").Value = "  -1.55%  "
$ws.Range("E30").Value = "  +0.65%  "
$ws.Range("E31").Value = "  -0.96%  "
$ws.Range("E32").Value = "  -1.56%  "
Set-TextValue $ws.Range("D33") "1.291.84"
$ws.Range("E33").Value = "  +1.35%  "
$ws.Range("E34").Value = "  -0.51%  "
$ws.Range("E35").Value = "  +1.12%  "
$ws.Range("E36").Value = "  -3.05%  "
Set-TextValue $ws.Range("D37") "0.534"
$ws.Range("E37").Value = "  +0.55%  "
Set-TextValue $ws.Range("D38") "0.816"
$ws.Range("E38").Value = "  -1.62%  "
$ws.Range("E39").Value = "  +0.09%  "
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws.Range("D40") "2.24"
$ws.Range("E40").Value = "  -0.08%  "
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue $ws.Range("D41") "0.802"
$ws.Range("E41").Value = "  -1.44%  "
$ws.Range("E42").Value = "  -2.50%  "
Set-TextValue $ws.Range("D43") "1.788.58"
$ws.Range("E43").Value = "  +0.36%  "
Set-TextValue $ws.Range("D44") "61.34"
$ws.Range("E44").Value = "  +3.23%  "
Set-TextValue $ws.Range("D45") "91.28"
$ws.Range("E46").Value = "  +0.02%  "
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D47") "0.0527"
$ws.Range("E47").Value = "  +2.21%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue $ws.Range("D48") "0.0₆0103"
$ws.Range("E48").Value = "  -2.33%  "
Set-TextValue $ws.Range("D49") "7.67"
$ws.Range("E49").Value = "  -1.61%  "
$ws.Range("E50").Value = "  -0.05%  "
$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue $ws.Range("D51") "0.406"
$ws.Range("E51").Value = "  -0.13%  "
